$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new tracker row (row 52) for goal G5 "Investment Plan"
$ws.Range("A52").Value = "G5"
$ws.Range("B52").Value = "Investment Plan"
$ws.Range("C52").Value = 45907
$ws.Range("C52").NumberFormat = "YYYY-MM-DD"
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 0
